$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("toDate")
$ws.Range("B9").Value = "hello"
Write-Host "done"
